# edit.ps1
# Applies the "New crime data collected" update to the CompStat_1 sheet:
#  - bumps the report Volume/Number and the covered-week dates in the header
#  - rewrites the Crime Complaints table (rows 14-30, cols C-N) with the
#    refreshed weekly figures, including cells that flip between a numeric
#    value and the "N/A" (0) / "***.*" text placeholders used elsewhere in
#    the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header: Volume 30 Number 4 -> 5 ; week 1/23/2023-1/29/2023 -> 1/30/2023-2/5/2023
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  5"
$ws.Range("C9").Value = "Report Covering the Week  1/30/2023  Through  2/5/2023"

# ---------------------------------------------------------------------
# Helper: turn a numeric cell into the shared "N/A"/"***.*" text marker
# while keeping the same right-aligned text style (style id 14) used by
# every other placeholder cell in the table (sourced from G14).
# ---------------------------------------------------------------------
function Set-Placeholder($ref, $text) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $text
    $ws.Range("G14").Copy() | Out-Null
    $c.PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------
# Helper: turn a text placeholder cell back into a genuine number while
# restoring the numeric style used by its row-mates ($styleSource).
# ---------------------------------------------------------------------
function Set-NumberFromPlaceholder($ref, $value, $styleSource) {
    $c = $ws.Range($ref)
    $c.Value = $value
    $ws.Range($styleSource).Copy() | Out-Null
    $c.PasteSpecial(-4122) | Out-Null
}

# Row 14 - Murder
Set-Placeholder "F14" "0"

# Row 15 - Rape
Set-Placeholder "C15" "0"

# Row 16 - Robbery
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 300
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = 100
$ws.Range("I16").Value = 20
$ws.Range("J16").Value = 10
$ws.Range("K16").Value = 100
$ws.Range("L16").Value = 150
$ws.Range("M16").Value = -31.034482758620
$ws.Range("N16").Value = -84.848484848484

# Row 17 - Fel. Assault
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 16
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = 6.666666666666
$ws.Range("I17").Value = 22
$ws.Range("J17").Value = 20
$ws.Range("K17").Value = 10
$ws.Range("L17").Value = 100
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = -54.166666666666

# Row 18 - Burglary
$ws.Range("C18").Value = 3
$ws.Range("E18").Value = -25
$ws.Range("F18").Value = 17
$ws.Range("G18").Value = 25
$ws.Range("H18").Value = -32
$ws.Range("I18").Value = 28
$ws.Range("J18").Value = 38
$ws.Range("K18").Value = -26.315789473684
$ws.Range("L18").Value = 40
$ws.Range("M18").Value = -31.707317073170
$ws.Range("N18").Value = -74.774774774774

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 17
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 70
$ws.Range("F19").Value = 71
$ws.Range("G19").Value = 42
$ws.Range("H19").Value = 69.047619047619
$ws.Range("I19").Value = 84
$ws.Range("J19").Value = 49
$ws.Range("K19").Value = 71.428571428571
$ws.Range("L19").Value = 147.058823529412
$ws.Range("M19").Value = 110
$ws.Range("N19").Value = 90.909090909090

# Row 20 - G.L.A.
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = -10
$ws.Range("I20").Value = 13
$ws.Range("J20").Value = 14
$ws.Range("K20").Value = -7.142857142857
$ws.Range("L20").Value = 8.333333333333
$ws.Range("M20").Value = -18.75
$ws.Range("N20").Value = -82.432432432432

# Row 21 - TOTAL (bold row, styles 18/19, values only)
$ws.Range("C21").Value = 30
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = 42.857142857142
$ws.Range("F21").Value = 130
$ws.Range("G21").Value = 101
$ws.Range("H21").Value = 28.712871287128
$ws.Range("I21").Value = 171
$ws.Range("J21").Value = 133
$ws.Range("K21").Value = 28.571428571428
$ws.Range("L21").Value = 94.318181818181
$ws.Range("M21").Value = 23.913043478260
$ws.Range("N21").Value = -58.894230769230

# Row 22 - Transit
Set-Placeholder "D22" "0"
Set-Placeholder "E22" "***.*"
$ws.Range("F22").Value = 1
$ws.Range("H22").Value = -75

# Row 23 - Housing
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = -66.666666666666
$ws.Range("F23").Value = 9
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 13
$ws.Range("J23").Value = 13
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 18.181818181818
$ws.Range("M23").Value = 8.333333333333

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 20
$ws.Range("E24").Value = -5
$ws.Range("F24").Value = 91
$ws.Range("H24").Value = 4.597701149425
$ws.Range("I24").Value = 115
$ws.Range("J24").Value = 107
$ws.Range("K24").Value = 7.476635514018
$ws.Range("L24").Value = 51.315789473684
$ws.Range("M24").Value = 0

# Row 25 - Misd. Assault
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 42.857142857142
$ws.Range("F25").Value = 31
$ws.Range("G25").Value = 34
$ws.Range("H25").Value = -8.823529411764
$ws.Range("I25").Value = 42
$ws.Range("J25").Value = 39
$ws.Range("K25").Value = 7.692307692307
$ws.Range("L25").Value = 23.529411764705
$ws.Range("M25").Value = 13.513513513513

# Row 26 - UCR Rape*
$ws.Range("C26").Value = 2
$ws.Range("F26").Value = 6
$ws.Range("H26").Value = 200
$ws.Range("I26").Value = 6
$ws.Range("K26").Value = 200
$ws.Range("L26").Value = 50

# Row 27 - Other Sex Crimes
Set-Placeholder "C27" "0"
Set-Placeholder "D27" "0"
Set-Placeholder "E27" "***.*"
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 50
$ws.Range("L27").Value = -40

# Row 28 - Shooting Vic.
Set-Placeholder "D28" "0"
Set-Placeholder "E28" "***.*"

# Row 29 - Shooting Inc.
Set-Placeholder "D29" "0"
Set-Placeholder "E29" "***.*"

# Row 30 - Hate Crimes (placeholders become real numbers this week)
Set-NumberFromPlaceholder "D30" 1 "G30"
Set-NumberFromPlaceholder "E30" -100 "H30"
$ws.Range("J30").Value = 2
